$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Il33"
$ws.Cells.Item(2, 3).Value = "Il1rl1"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.4785103333333333
$ws.Cells.Item(2, 8).Value = 1.435531
$ws.Cells.Item(2, 9).Value = 0.01183374253319748
$ws.Cells.Item(2, 10).Value = 0.01183374253319748
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 1.557338
$ws.Cells.Item(2, 14).Value = 4.672014
$ws.Cells.Item(2, 15).Value = 0.3240096237227595
$ws.Cells.Item(2, 16).Value = 0.3240096237227596
$ws.Cells.Item(2, 17).Value = 0.7452023254926666
$ws.Cells.Item(2, 18).Value = 6.706820929434
$ws.Cells.Item(2, 19).Value = 0.003834246465413332
$ws.Cells.Item(2, 20).Value = 0.003834246465413333

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Il33"
$ws.Cells.Item(3, 3).Value = "Il1rl1"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.4785103333333333
$ws.Cells.Item(3, 8).Value = 1.435531
$ws.Cells.Item(3, 9).Value = 0.01183374253319748
$ws.Cells.Item(3, 10).Value = 0.01183374253319748
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.249118
$ws.Cells.Item(3, 14).Value = 9.747354
$ws.Cells.Item(3, 15).Value = 0.6759903762772405
$ws.Cells.Item(3, 16).Value = 0.6759903762772405
$ws.Cells.Item(3, 17).Value = 1.554736537219333
$ws.Cells.Item(3, 18).Value = 13.992628834974
$ws.Cells.Item(3, 19).Value = 0.007999496067784152
$ws.Cells.Item(3, 20).Value = 0.007999496067784152

# Row 4
$ws.Cells.Item(4, 1).Value = "sCs"
$ws.Cells.Item(4, 2).Value = "Il33"
$ws.Cells.Item(4, 3).Value = "Il1rl1"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 37.79047733333334
$ws.Cells.Item(4, 8).Value = 113.371432
$ws.Cells.Item(4, 9).Value = 0.9345728771499232
$ws.Cells.Item(4, 10).Value = 0.9345728771499232
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 1.557338
$ws.Cells.Item(4, 14).Value = 4.672014
$ws.Cells.Item(4, 15).Value = 0.3240096237227595
$ws.Cells.Item(4, 16).Value = 0.3240096237227596
$ws.Cells.Item(4, 17).Value = 58.85254638933866
$ws.Cells.Item(4, 18).Value = 529.672917504048
$ws.Cells.Item(4, 19).Value = 0.3028106062668434
$ws.Cells.Item(4, 20).Value = 0.3028106062668435

# Row 5
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Il33"
$ws.Cells.Item(5, 3).Value = "Il1rl1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 37.79047733333334
$ws.Cells.Item(5, 8).Value = 113.371432
$ws.Cells.Item(5, 9).Value = 0.9345728771499232
$ws.Cells.Item(5, 10).Value = 0.9345728771499232
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 3.249118
$ws.Cells.Item(5, 14).Value = 9.747354
$ws.Cells.Item(5, 15).Value = 0.6759903762772405
$ws.Cells.Item(5, 16).Value = 0.6759903762772405
$ws.Cells.Item(5, 17).Value = 122.7857201323253
$ws.Cells.Item(5, 18).Value = 1105.071481190928
$ws.Cells.Item(5, 19).Value = 0.6317622708830798
$ws.Cells.Item(5, 20).Value = 0.6317622708830798

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Il33"
$ws.Cells.Item(6, 3).Value = "Il1rl1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.167107
$ws.Cells.Item(6, 8).Value = 6.501321
$ws.Cells.Item(6, 9).Value = 0.05359338031687927
$ws.Cells.Item(6, 10).Value = 0.05359338031687926
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 1.557338
$ws.Cells.Item(6, 14).Value = 4.672014
$ws.Cells.Item(6, 15).Value = 0.3240096237227595
$ws.Cells.Item(6, 16).Value = 0.3240096237227596
$ws.Cells.Item(6, 17).Value = 3.374918081166
$ws.Cells.Item(6, 18).Value = 30.374262730494
$ws.Cells.Item(6, 19).Value = 0.0173647709905028
$ws.Cells.Item(6, 20).Value = 0.0173647709905028

# Row 7
$ws.Cells.Item(7, 1).Value = "ECs"
$ws.Cells.Item(7, 2).Value = "Il33"
$ws.Cells.Item(7, 3).Value = "Il1rl1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.167107
$ws.Cells.Item(7, 8).Value = 6.501321
$ws.Cells.Item(7, 9).Value = 0.05359338031687927
$ws.Cells.Item(7, 10).Value = 0.05359338031687926
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.249118
$ws.Cells.Item(7, 14).Value = 9.747354
$ws.Cells.Item(7, 15).Value = 0.6759903762772405
$ws.Cells.Item(7, 16).Value = 0.6759903762772405
$ws.Cells.Item(7, 17).Value = 7.041186361626
$ws.Cells.Item(7, 18).Value = 63.370677254634
$ws.Cells.Item(7, 19).Value = 0.03622860932637647
$ws.Cells.Item(7, 20).Value = 0.03622860932637646
